# This script applies a stock-count correction to the "CryCompanywiseStockReport"
# workbook. A physical stock recount changed the on-hand Quantity (column F) for a
# number of SKUs; the corresponding stock Value (column G = Rate(D) x Qty(F)) is
# recalculated for each affected line, each company "Sub Total:" (column B) is
# recalculated as the sum of the Value column for that company's block, and the
# overall "Sub Total:"/"Grand Total:" rows (B1099/B1100) are recalculated as the
# sum of every company sub-total. Two rows (228/229) had their Item Code (column B)
# values swapped to correct a data-entry mix-up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F45").Value2 = 125
$ws.Range("G45").Value2 = 3903.75
$ws.Range("F46").Value2 = 157
$ws.Range("G46").Value2 = 10996.28
$ws.Range("F63").Value2 = 206
$ws.Range("G63").Value2 = 19269.24
$ws.Range("B77").Value2 = 100658.76
$ws.Range("F83").Value2 = 5
$ws.Range("G83").Value2 = 4269.85
$ws.Range("B95").Value2 = 34693.47
$ws.Range("F114").Value2 = 10
$ws.Range("G114").Value2 = 2528.4
$ws.Range("B159").Value2 = 459430.35
$ws.Range("F183").Value2 = 2
$ws.Range("G183").Value2 = 202.46
$ws.Range("B195").Value2 = 35643.54
$ws.Range("F205").Value2 = 124
$ws.Range("G205").Value2 = 6135.52
$ws.Range("F206").Value2 = 58
$ws.Range("G206").Value2 = 2869.84
$ws.Range("F207").Value2 = 88
$ws.Range("G207").Value2 = 3724.16
$ws.Range("F214").Value2 = 22
$ws.Range("G214").Value2 = 2014.54
$ws.Range("B217").Value2 = 45795
$ws.Range("B228").Value2 = 57756
$ws.Range("B229").Value2 = 53925
$ws.Range("F240").Value2 = 12
$ws.Range("G240").Value2 = 4789.92
$ws.Range("B256").Value2 = 251069.69
$ws.Range("F295").Value2 = 14
$ws.Range("G295").Value2 = 650.86
$ws.Range("B305").Value2 = 7256.35
$ws.Range("F316").Value2 = 160
$ws.Range("G316").Value2 = 12355.2
$ws.Range("F317").Value2 = 173
$ws.Range("G317").Value2 = 11921.43
$ws.Range("F318").Value2 = 202
$ws.Range("G318").Value2 = 10186.86
$ws.Range("B320").Value2 = 51923.84
$ws.Range("F332").Value2 = 6
$ws.Range("G332").Value2 = 4852.08
$ws.Range("B336").Value2 = 66237.98
$ws.Range("F360").Value2 = 23
$ws.Range("G360").Value2 = 3117.42
$ws.Range("F364").Value2 = 10
$ws.Range("G364").Value2 = 820.3
$ws.Range("F380").Value2 = 71
$ws.Range("G380").Value2 = 1434.91
$ws.Range("B409").Value2 = 201897.77
$ws.Range("F413").Value2 = 501
$ws.Range("G413").Value2 = 51457.71
$ws.Range("F418").Value2 = 160
$ws.Range("G418").Value2 = 13270.4
$ws.Range("F423").Value2 = 96
$ws.Range("G423").Value2 = 15232.32
$ws.Range("F427").Value2 = 198
$ws.Range("G427").Value2 = 7896.24
$ws.Range("F431").Value2 = 53
$ws.Range("G431").Value2 = 5980.52
$ws.Range("F449").Value2 = 46
$ws.Range("G449").Value2 = 5051.72
$ws.Range("F455").Value2 = 4
$ws.Range("G455").Value2 = 497.08
$ws.Range("F458").Value2 = 180
$ws.Range("G458").Value2 = 15429.6
$ws.Range("F462").Value2 = 16
$ws.Range("G462").Value2 = 4222.88
$ws.Range("F471").Value2 = 158
$ws.Range("G471").Value2 = 16417.78
$ws.Range("F472").Value2 = 181
$ws.Range("G472").Value2 = 21472.03
$ws.Range("F473").Value2 = 439
$ws.Range("G473").Value2 = 25958.07
$ws.Range("F475").Value2 = 60
$ws.Range("G475").Value2 = 29662.8
$ws.Range("F476").Value2 = 1827
$ws.Range("G476").Value2 = 39627.63
$ws.Range("F477").Value2 = 1780
$ws.Range("G477").Value2 = 10697.8
$ws.Range("F478").Value2 = 391
$ws.Range("G478").Value2 = 32492.1
$ws.Range("F480").Value2 = 89
$ws.Range("G480").Value2 = 37638.1
$ws.Range("F483").Value2 = 375
$ws.Range("G483").Value2 = 26366.25
$ws.Range("F484").Value2 = 169
$ws.Range("G484").Value2 = 28954.77
$ws.Range("F485").Value2 = 275
$ws.Range("G485").Value2 = 41571.75
$ws.Range("F486").Value2 = 113
$ws.Range("G486").Value2 = 15418.85
$ws.Range("F490").Value2 = 173
$ws.Range("G490").Value2 = 16320.82
$ws.Range("F491").Value2 = 509
$ws.Range("G491").Value2 = 10342.88
$ws.Range("F494").Value2 = 173
$ws.Range("G494").Value2 = 6876.75
$ws.Range("B496").Value2 = 1132480.94
$ws.Range("F509").Value2 = 36
$ws.Range("G509").Value2 = 12939.48
$ws.Range("F510").Value2 = 16
$ws.Range("G510").Value2 = 2596.64
$ws.Range("F512").Value2 = 53
$ws.Range("G512").Value2 = 11566.19
$ws.Range("B513").Value2 = 105681.83
$ws.Range("F556").Value2 = 3
$ws.Range("G556").Value2 = 1379.97
$ws.Range("B566").Value2 = 23187.23
$ws.Range("F568").Value2 = 127
$ws.Range("G568").Value2 = 8718.549999999999
$ws.Range("F570").Value2 = 65
$ws.Range("G570").Value2 = 6860.1
$ws.Range("F573").Value2 = 46
$ws.Range("G573").Value2 = 1256.72
$ws.Range("F578").Value2 = 466
$ws.Range("G578").Value2 = 45015.6
$ws.Range("B584").Value2 = 105652.7
$ws.Range("F639").Value2 = 478
$ws.Range("G639").Value2 = 12571.4
$ws.Range("B654").Value2 = 143290.65
$ws.Range("F656").Value2 = 152
$ws.Range("G656").Value2 = 5619.44
$ws.Range("B673").Value2 = 21653.99
$ws.Range("F680").Value2 = 136
$ws.Range("G680").Value2 = 5831.68
$ws.Range("B697").Value2 = 92182.34
$ws.Range("F773").Value2 = 164
$ws.Range("G773").Value2 = 588.76
$ws.Range("F774").Value2 = 63
$ws.Range("G774").Value2 = 904.05
$ws.Range("F776").Value2 = 270
$ws.Range("G776").Value2 = 5969.7
$ws.Range("F777").Value2 = 140
$ws.Range("G777").Value2 = 1003.8
$ws.Range("F784").Value2 = 47
$ws.Range("G784").Value2 = 1213.54
$ws.Range("F791").Value2 = 31
$ws.Range("G791").Value2 = 3255.62
$ws.Range("B792").Value2 = 32918.35
$ws.Range("F796").Value2 = 90
$ws.Range("G796").Value2 = 25728.3
$ws.Range("F811").Value2 = 184
$ws.Range("G811").Value2 = 13533.2
$ws.Range("F815").Value2 = 129
$ws.Range("G815").Value2 = 15845.07
$ws.Range("F817").Value2 = 47
$ws.Range("G817").Value2 = 3901.94
$ws.Range("B818").Value2 = 163094.22
$ws.Range("F943").Value2 = 7
$ws.Range("G943").Value2 = 11237.1
$ws.Range("F944").Value2 = 5
$ws.Range("G944").Value2 = 6232.65
$ws.Range("F945").Value2 = 0
$ws.Range("G945").Value2 = 0
$ws.Range("F948").Value2 = 47
$ws.Range("G948").Value2 = 6141.49
$ws.Range("F950").Value2 = 52
$ws.Range("G950").Value2 = 12764.44
$ws.Range("F951").Value2 = 2
$ws.Range("G951").Value2 = 3201.8
$ws.Range("B954").Value2 = 137942.22
$ws.Range("F956").Value2 = 16
$ws.Range("G956").Value2 = 1304.96
$ws.Range("F967").Value2 = 285
$ws.Range("G967").Value2 = 6190.2
$ws.Range("F972").Value2 = 216
$ws.Range("G972").Value2 = 8953.200000000001
$ws.Range("F978").Value2 = 688
$ws.Range("G978").Value2 = 99072
$ws.Range("B982").Value2 = 297740.48
$ws.Range("F990").Value2 = 29
$ws.Range("G990").Value2 = 2448.76
$ws.Range("F1015").Value2 = 60
$ws.Range("G1015").Value2 = 3419.4
$ws.Range("B1016").Value2 = 142048.86
$ws.Range("B1099").Value2 = 6371169.67
$ws.Range("B1100").Value2 = 6371169.67
